$wb = $excel.ActiveWorkbook

# --- Nedas_J: add Week 5 row ---
$ws = $wb.Worksheets.Item("Nedas_J")
$ws.Range("A5").Value = "Week 5"
$ws.Range("B5").Value = 0.58333333333333337
$ws.Range("B5").NumberFormat = "h:mm"

# --- Adomas_J: add Week 5 row ---
$ws = $wb.Worksheets.Item("Adomas_J")
$ws.Range("A5").Value = "Week 5"
$ws.Range("B5").Value = 0.30555555555555552
$ws.Range("B5").NumberFormat = "h:mm"

# --- Aiste_G: add Week 5 row (3rd sheet; name has a diacritic, use index to be safe) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("A5").Value = "Week 5"
$ws.Range("B5").Value = 0.15277777777777776
$ws.Range("B5").NumberFormat = "h:mm"

# --- Gabrielius_D: add Week 5 row ---
$ws = $wb.Worksheets.Item("Gabrielius_D")
$ws.Range("A5").Value = "Week 5"
$ws.Range("B5").Value = 0.20833333333333334
$ws.Range("B5").NumberFormat = "h:mm"

# --- Overview: add Total column + Week 5 rows ---
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("D1").Value = "Total"
$ws.Range("E1").Formula = "=SUM(B2+B4+B6+B8+B10)"
$ws.Range("E1").NumberFormat = "[hh]:mm"

$ws.Range("A9").Value = "Week 5 (working in a group)"
$ws.Range("B9").Value = 0.0625
$ws.Range("B9").NumberFormat = "h:mm"

$ws.Range("A10").Value = "Week 5 (total working hours)"
$ws.Range("B8").Copy($ws.Range("B10"))
$nameNedas = $wb.Worksheets.Item(1).Name
$nameAdomas = $wb.Worksheets.Item(2).Name
$nameAiste = $wb.Worksheets.Item(3).Name
$nameGabrielius = $wb.Worksheets.Item(4).Name
$ws.Range("B10").Formula = "=SUM(" + $nameNedas + "!B5, " + $nameAdomas + "!B5, " + $nameAiste + "!B5, " + $nameGabrielius + "!B5) + 4 * B9"

$wb.Save()
